# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the aggregated "全部类型" sheet, per the scraped-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 24
$ws1.Range("F5").Value = 3271
$ws1.Range("F9").Value = 25
$ws1.Range("F10").Value = 1201
$ws1.Range("F12").Value = 1187

# --- Sheet "全部类型" (all types, aggregate of every sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 24
$ws4.Range("F5").Value = 3271
$ws4.Range("F10").Value = 25
$ws4.Range("F11").Value = 1201
$ws4.Range("F13").Value = 1187
